$wb = $excel.ActiveWorkbook

$oldVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: " + $newVersion

$about.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Tianan No. 1 Coal Mine, China, M1208, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 7; $r++) {
    $data.Range("S$r").Value = $newVersion
}
